$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 8334012
$ws.Range("I12").Value = 14285894
$ws.Range("J12").Value = 1377.6
$ws.Range("K12").Value = 14285894
$ws.Range("L12").Value = 1377.6
$ws.Range("M12").Value = -14285724
$ws.Range("N12").Value = -1717.6

$ws.Range("H86").Value = 3011.7856
$ws.Range("I86").Value = 2379
$ws.Range("K86").Value = 2379
$ws.Range("M86").Value = -1256

$ws.Range("H89").Value = 3011.7856
$ws.Range("I89").Value = 2379
$ws.Range("K89").Value = 11895
$ws.Range("M89").Value = -6279

$ws.Range("H106").Value = 2112
$ws.Range("I106").Value = 2106.5
$ws.Range("K106").Value = 2106.5
$ws.Range("M106").Value = -1475.5

$ws.Range("H111").Value = 3741.8333
$ws.Range("J111").Value = 4030.6
$ws.Range("L111").Value = 12091.8
$ws.Range("N111").Value = -18225.8

$ws.Range("H112").Value = 3671.4443
$ws.Range("J112").Value = 3630.375
$ws.Range("L112").Value = 10891.125
$ws.Range("N112").Value = -13107.125

$ws.Range("H137").Value = 1567
$ws.Range("I137").Value = 1445.5294
$ws.Range("K137").Value = 4336.5882
$ws.Range("M137").Value = -1786.5882

$ws.Range("H138").Value = 2154
$ws.Range("I138").Value = 1576.1428
$ws.Range("J138").Value = 2378.7222
$ws.Range("K138").Value = 4728.428400000001
$ws.Range("L138").Value = 7136.1666
$ws.Range("M138").Value = 411.5715999999993
$ws.Range("N138").Value = -17416.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 49999.332
$ws.Range("J24").Value = 49999.332
$ws.Range("L24").Value = 49999.332
$ws.Range("N24").Value = -50747.332

$ws.Range("H32").Value = 2015.2931
$ws.Range("I32").Value = 2057.0178
$ws.Range("K32").Value = 2057.0178
$ws.Range("M32").Value = -1770.0178

$ws.Range("I45").Value = 2114.9
$ws.Range("J45").Value = 7937.5
$ws.Range("K45").Value = 2114.9
$ws.Range("L45").Value = 7937.5
$ws.Range("M45").Value = -1737.9
$ws.Range("N45").Value = -8691.5

$ws.Range("H74").Value = 3812.55
$ws.Range("I74").Value = 2946.0417
$ws.Range("J74").Value = 5112.3125
$ws.Range("K74").Value = 2946.0417
$ws.Range("L74").Value = 5112.3125
$ws.Range("M74").Value = -2072.0417
$ws.Range("N74").Value = -6860.3125

$ws.Range("H77").Value = 3812.55
$ws.Range("I77").Value = 2946.0417
$ws.Range("J77").Value = 5112.3125
$ws.Range("K77").Value = 14730.2085
$ws.Range("L77").Value = 25561.5625
$ws.Range("M77").Value = -10362.2085
$ws.Range("N77").Value = -34297.5625

$ws.Range("H94").Value = 59999
$ws.Range("J94").Value = 59999
$ws.Range("L94").Value = 59999
$ws.Range("N94").Value = -61801

$ws.Range("H96").Value = 40999.668
$ws.Range("J96").Value = 40999.668
$ws.Range("L96").Value = 40999.668
$ws.Range("N96").Value = -46491.668

$ws.Range("H100").Value = 49999.332
$ws.Range("J100").Value = 49999.332
$ws.Range("L100").Value = 49999.332
$ws.Range("N100").Value = -52163.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1871.4286
$ws.Range("I86").Value = 1871.4286
$ws.Range("K86").Value = 1871.4286
$ws.Range("M86").Value = -748.4286

$ws.Range("H89").Value = 1871.4286
$ws.Range("I89").Value = 1871.4286
$ws.Range("K89").Value = 9357.143
$ws.Range("M89").Value = -3741.143

$ws.Range("H99").Value = 5457.6313
$ws.Range("I99").Value = 4108.636
$ws.Range("K99").Value = 4108.636
$ws.Range("M99").Value = -2610.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3708.5454
$ws.Range("I31").Value = 2732.6191
$ws.Range("J31").Value = 5416.4165
$ws.Range("K31").Value = 2732.6191
$ws.Range("L31").Value = 5416.4165
$ws.Range("M31").Value = -2437.6191
$ws.Range("N31").Value = -6006.4165

$ws.Range("H34").Value = 3708.5454
$ws.Range("I34").Value = 2732.6191
$ws.Range("J34").Value = 5416.4165
$ws.Range("K34").Value = 2732.6191
$ws.Range("L34").Value = 5416.4165
$ws.Range("M34").Value = -2530.6191
$ws.Range("N34").Value = -5820.4165

$ws.Range("H58").Value = 6931.4
$ws.Range("I58").Value = 4840.2354
$ws.Range("K58").Value = 4840.2354
$ws.Range("M58").Value = -4637.2354

$ws.Range("H136").Value = 6931.4
$ws.Range("I136").Value = 4840.2354
$ws.Range("K136").Value = 14520.7062
$ws.Range("M136").Value = -11970.7062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 260
$ws.Range("J12").Value = 454.25
$ws.Range("L12").Value = 1362.75
$ws.Range("N12").Value = -1708.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 54000
$ws.Range("J34").Value = 54000
$ws.Range("L34").Value = 54000
$ws.Range("N34").Value = -54536

$ws.Range("H69").Value = 33995
$ws.Range("J69").Value = 33995
$ws.Range("L69").Value = 33995
$ws.Range("N69").Value = -35493

$ws.Range("H72").Value = 33995
$ws.Range("J72").Value = 33995
$ws.Range("L72").Value = 101985
$ws.Range("N72").Value = -109473

$ws.Range("H76").Value = 54000
$ws.Range("J76").Value = 54000
$ws.Range("L76").Value = 54000
$ws.Range("N76").Value = -54630

$ws.Range("H79").Value = 54000
$ws.Range("J79").Value = 54000
$ws.Range("L79").Value = 54000
$ws.Range("N79").Value = -56184

$ws.Range("H132").Value = 1567.725
$ws.Range("I132").Value = 1567.725
$ws.Range("K132").Value = 4703.174999999999
$ws.Range("M132").Value = -2173.174999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4362.9614
$ws.Range("I136").Value = 3933.818
$ws.Range("K136").Value = 11801.454
$ws.Range("M136").Value = -9251.454000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 785.3333
$ws.Range("I107").Value = 450.75
$ws.Range("K107").Value = 1352.25
$ws.Range("M107").Value = 567.75

$ws.Range("H113").Value = 558.8333
$ws.Range("I113").Value = 577.26666
$ws.Range("K113").Value = 1731.79998
$ws.Range("M113").Value = 438.20002

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 5080.75
$ws.Range("I122").Value = 4794.8
$ws.Range("K122").Value = 14384.4
$ws.Range("M122").Value = -11934.4

$ws.Range("H136").Value = 4112.4136
$ws.Range("I136").Value = 2323.4583
$ws.Range("J136").Value = 12699.4
$ws.Range("K136").Value = 6970.374899999999
$ws.Range("L136").Value = 38098.2
$ws.Range("M136").Value = -4420.374899999999
$ws.Range("N136").Value = -43198.2
Write-Output "Applied all changes"
